$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# The title text is currently split across multiple runs ("Two-Column", " ",
# "Layout") that together read as "Two-Column Layout". Assigning that exact
# same string back as a no-op wouldn't force the runs to merge, so first
# assign a distinct placeholder to force a real text rewrite (collapsing to
# a single run), then assign the final desired text.
$tr.Text = "placeholder"
$tr.Text = "Two-Column Layout"
